$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert two new bullet points right after "Enter assumed true population
#    values..." (paragraph 8), before the existing "F should be zero" bullet.
#    InsertParagraphAfter() on a ListParagraph-styled range inherits the
#    pStyle/numPr automatically, so the new paragraphs come out as proper
#    list items.
# ---------------------------------------------------------------------------
$pEnter = $d.Paragraphs.Item(8)
$pEnter.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(9)
$pNew1.Range.Text = "Linf, K, (t0), Amax from the Life History of Bottomfish Management Unit Species of Guam. NMFS-PIFSC-170. CV assumed 2.5% for Linf, 5.7% for K, 25% for L0."

$pNew1b = $d.Paragraphs.Item(9)
$pNew1b.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item(10)
$pNew2.Range.Text = "M calculated from Amax following Then et al. 2015 M = 4.899*Amax^-0.916"

# ---------------------------------------------------------------------------
# After the two insertions above, paragraphs shift down by 2:
#   11 = "F should be zero"                (was 9, unchanged text)
#   12 = "mincat, catsd, ..."              (was 10, unchanged text)
#   13 = "Reminder: age_max ..."           (was 11, unchanged text)
#   14 = "N can be as large ..."           (was 12, needs simplifying)
#   15 = "Linf_k_cor_TF is a TRUE ..."     (was 13, needs mid insertion)
#   16 = empty                             (was 14)
#   17 = "2. Save each simulated ..."      (was 15, needs reword)
#   18 = empty                             (was 16)
#   19 = empty                             (was 17)
#   20 = "To do: fill in population ..."   (was 18)
# ---------------------------------------------------------------------------

# 2. Simplify the "N can be as large..." run-time bullet down to one run.
$pRunTime = $d.Paragraphs.Item(14)
$rRunTime = $pRunTime.Range
$okRunTime = $rRunTime.Find.Execute("N can be as large as possible, but will increase run time. for Amax = 55, each 100k takes about 3 minutes to run. ", $false, $false, $false, $false, $false, $true, 1, $false, "N can be as large as possible, but will increase run time. for Amax = 55, each 100k takes about 3 minutes to run. ", 2)

# 3. Insert the Jensen 1997 explanation in the middle of the Linf_k_cor_TF bullet.
$pLinf = $d.Paragraphs.Item(15)
$rLinf = $pLinf.Range
$jensenInsert = ". If Linf_k_cor_TF = TRUE then the Jensen 1997 Eq. 4 theoretical expected relationship between Linf and K (Linf = C1*k^(-1/LW_beta)) is used where LW_beta is the theoretical allometric W-L scalar, assume = 3 here and C1 is calculated from the mean population values of Linf and k"
$okLinf = $rLinf.Find.Execute("(true). If Linf_k_cor_TF = TRUE,", $false, $false, $false, $false, $false, $true, 1, $false, "(true)$jensenInsert. If Linf_k_cor_TF = TRUE,", 2)

# 4. Reword the "2. Save each simulated population dataframe..." bullet.
$pSave = $d.Paragraphs.Item(17)
$rSave = $pSave.Range
$okSave = $rSave.Find.Execute("each simulated population dataframe in its own workspace.", $false, $false, $false, $false, $false, $true, 1, $false, "all simulated population dataframes in a workspace.", 2)

# ---------------------------------------------------------------------------
# 5. After "2. Save..." (paragraph 17), insert a blank paragraph followed by
#    a new "To do" paragraph with the expanded text.
# ---------------------------------------------------------------------------
$pSave.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item(18)
$pBlank.Range.InsertParagraphAfter()
$pToDoNew = $d.Paragraphs.Item(19)
$pToDoNew.Range.Text = "To do: fill in population parameter assumptions, make a pdf with figures summarizing population assumptions."

# ---------------------------------------------------------------------------
# After the insertions above, the two originally-empty trailing paragraphs
# and the final "To do: fill in population parameter assumptions." paragraph
# shift down by 2:
#   20 = empty                              (was 16)
#   21 = empty                              (was 17)
#   22 = "To do: fill in population ..."    (was 18, becomes "3. run with every update:")
# ---------------------------------------------------------------------------

# 6. Replace the final paragraph's text with the new "3. run with every update:" heading.
$pRunUpdate = $d.Paragraphs.Item(22)
$rRunUpdate = $pRunUpdate.Range
$okRunUpdate = $rRunUpdate.Find.Execute("To do: fill in population parameter assumptions.", $false, $false, $false, $false, $false, $true, 1, $false, "3. run with every update:", 2)

# 7. Append two new closing paragraphs after it.
$pRunUpdate.Range.InsertParagraphAfter()
$pColor = $d.Paragraphs.Item(23)
$pColor.Range.Text = "Add a color code to show added samples in the hist during survey"

$pColor2 = $d.Paragraphs.Item(23)
$pColor2.Range.InsertParagraphAfter()
$pFig = $d.Paragraphs.Item(24)
$pFig.Range.Text = "show fig with units, provide calculated relative error as a number"

Write-Output "paragraphs=$($d.Paragraphs.Count) runTime=$okRunTime linf=$okLinf save=$okSave runUpdate=$okRunUpdate"
